# Updated cryptos list on Sun Mar 31 19:52:05 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns with newer
# figures, and - because ThetaToken/ApeXProtocol swapped rank - also
# rewrites B/C/D/E for rows 45 and 46 so the two coins trade places.
#
# Price-column values are written with a leading "'" (quote-prefix) so
# numeric-looking text such as "1.00" or "7.20" stays plain text instead
# of being auto-coerced to a number (which would silently drop the
# formatted trailing zeros, e.g. "1.00" -> 1, "7.20" -> 7.2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''70.973.18'
$ws.Range("E2").Value = '  +1.90%  '
$ws.Range("D3").Value = '''3.637.80'
$ws.Range("E3").Value = '  +3.92%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''605.16'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").Value = '''199.62'
$ws.Range("E6").Value = '  +2.18%  '
$ws.Range("D7").Value = '''0.629'
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '''0.223'
$ws.Range("E9").Value = '  +11.62%  '
$ws.Range("D10").Value = '''0.647'
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").Value = '''53.96'
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("E12").Value = '  +2.10%  '
$ws.Range("D13").Value = '''9.58'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").Value = '''4.207.74'
$ws.Range("E14").Value = '  +3.58%  '
$ws.Range("D15").Value = '''643.56'
$ws.Range("E15").Value = '  +8.44%  '
$ws.Range("D16").Value = '''13.01'
$ws.Range("E16").Value = '  +1.76%  '
$ws.Range("D17").Value = '''71.054.36'
$ws.Range("E17").Value = '  +1.75%  '
$ws.Range("D18").Value = '''3.665.44'
$ws.Range("E18").Value = '  +4.78%  '
$ws.Range("D19").Value = '''19.08'
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("E21").Value = '  +1.35%  '
$ws.Range("D22").Value = '''18.69'
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("D23").Value = '''5.36'
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("D24").Value = '''104.39'
$ws.Range("E24").Value = '  +2.16%  '
$ws.Range("D25").Value = '''4.64'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D27").Value = '''10.47'
$ws.Range("E27").Value = '  -3.39%  '
$ws.Range("D28").Value = '''9.77'
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("D29").Value = '''34.05'
$ws.Range("E29").Value = '  +2.37%  '
$ws.Range("D30").Value = '''4.77'
$ws.Range("E30").Value = '  +11.12%  '
$ws.Range("D31").Value = '''7.20'
$ws.Range("E31").Value = '  +2.52%  '
$ws.Range("D32").Value = '''12.25'
$ws.Range("E32").Value = '  -1.13%  '
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("D34").Value = '''63.39'
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("D35").Value = '''4.031.52'
$ws.Range("E35").Value = '  +8.58%  '
$ws.Range("E36").Value = '  +6.10%  '
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").Value = '''3.05'
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("D39").Value = '''511.09'
$ws.Range("E39").Value = '  +8.42%  '
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("D41").Value = '''36.77'
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("E44").Value = '  +2.00%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").Value = '''3.05'
$ws.Range("E45").Value = '  +8.58%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '''3.48'
$ws.Range("E46").Value = '  +6.11%  '
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").Value = '''8.68'
$ws.Range("E48").Value = '  +3.27%  '
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("D50").Value = '''0.000251'
$ws.Range("E50").Value = '  +2.68%  '
$ws.Range("E51").Value = '  +5.41%  '
